{"js": "// 1. Body: \"Overview of Phase 1b\" -> \"Overview of Phase 1a\"\nconst body = context.document.body;\nconst titleResults = body.search(\"Overview of Phase 1b\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nfor (const r of titleResults.items) {\n  r.insertText(\"Overview of Phase 1a\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. First-page footer: join \" \" and \"Pourassad Mohammadhossein\" into a\n// single contiguous run \" Pourassad Mohammadhossein\" (text itself is\n// unchanged, only the run split goes away).\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  const footer = section.getFooter(Word.HeaderFooterType.firstPage);\n  const nameResults = footer.search(\" Pourassad Mohammadhossein\", { matchCase: true });\n  nameResults.load(\"items\");\n  await context.sync();\n\n  for (const r of nameResults.items) {\n    r.insertText(\" Pourassad Mohammadhossein\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Body: \"Overview of Phase 1b\" -> \"Overview of Phase 1a\"\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"Overview of Phase 1b\")\nif ($r.Find.Found) {\n    $r.Text = \"Overview of Phase 1a\"\n}\n\n# 2. First-page footer: join \" \" and \"Pourassad Mohammadhossein\" into a\n# single contiguous run \" Pourassad Mohammadhossein\" (visible text is\n# unchanged, only the run split goes away).\nforeach ($section in $d.Sections) {\n    $footer = $section.Footers.Item(2)  # wdHeaderFooterFirstPage\n\n    $fr = $footer.Range\n    $fr.Find.ClearFormatting()\n    $fr.Find.Execute(\" Pourassad Mohammadhossein\")\n    if ($fr.Find.Found) {\n        # Assigning the exact same text back is a no-op, so the two runs\n        # would stay split. Go through a placeholder value first so the\n        # assignment actually changes the text and the runs get merged,\n        # then fix the placeholder back to the real name. Re-fetch the\n        # footer range before the second Find since $fr does not refresh\n        # its bounds after the first assignment.\n        $fr.Text = \" Pourassad Mohammadhossein#\"\n\n        $fr2 = $footer.Range\n        $fr2.Find.ClearFormatting()\n        $fr2.Find.Execute(\" Pourassad Mohammadhossein#\")\n        if ($fr2.Find.Found) {\n            $fr2.Text = \" Pourassad Mohammadhossein\"\n        }\n    }\n}\n"}
